$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Move the "Buy" block (old A6:D8 on Sheet1) onto Sheet2 first,
#     before we touch/clear the source cells, so the formatting/styles
#     travel with the values. ---
$ws1.Range("A6").Copy($ws2.Range("A1"))
$ws1.Range("B7").Copy($ws2.Range("B2"))
$ws1.Range("C7").Copy($ws2.Range("C2"))
$ws1.Range("D7").Copy($ws2.Range("D2"))
$ws1.Range("B8").Copy($ws2.Range("B3"))
$ws1.Range("C8").Copy($ws2.Range("C3"))
$ws1.Range("D8").Copy($ws2.Range("D3"))

# Column widths on the new Buy sheet (matches the bestFit widths used
# elsewhere in the workbook for this same data).
$ws2.Columns("B").ColumnWidth = 20.6
$ws2.Columns("C").ColumnWidth = 8.6

# --- Clean up Sheet1: drop the old rows 6-8 (now duplicated on Sheet2) ---
$ws1.Range("A6:D8").Clear()

# Sheet1 header text: "Login_App" -> "Login"
$ws1.Range("A1").Value = "Login"

# Fix the text on the moved header cell: "Buy" -> "buy"
$ws2.Range("A1").Value = "buy"

# New bottom-of-sheet note cell on Sheet1
$ws1.Range("C18").Value = "s"

# --- Rename the tabs ---
$ws1.Name = "Login_App"
$ws2.Name = "Buy"

# --- Selections, to match the saved view state ---
# (Select Sheet1's cell last so Sheet1 remains the active/tab-selected sheet)
[void]$ws2.Range("B18").Select()
[void]$ws1.Range("C18").Select()
